$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new test result values in column F for rows 14 and 15
$ws.Range("F14").Value = 0.0002281
$ws.Range("F15").Value = 0.3861

# Move the active selection to F15, matching the new selection state
$ws.Range("F15").Select()
